# Update "想去人数" (number wanting to attend) figures for several
# events in both the "展览" (Exhibition) and "全部类型" (All types)
# worksheets, per the upstream gh-pages data refresh at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 31
$ws1.Range("F5").Value = 3331
$ws1.Range("F7").Value = 406
$ws1.Range("F9").Value = 37
$ws1.Range("F10").Value = 23
$ws1.Range("F11").Value = 1231
$ws1.Range("F13").Value = 1331

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 31
$ws4.Range("F5").Value = 3331
$ws4.Range("F7").Value = 406
$ws4.Range("F10").Value = 37
$ws4.Range("F11").Value = 23
$ws4.Range("F14").Value = 1231
$ws4.Range("F16").Value = 1331
